$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report gains two newer daily snapshot columns ("Jun_17" and "Jun_15").
# They are inserted immediately before the existing first data column
# ("Jun_13" in column B), which — together with the column after it
# ("Jun_10" in column C) — shifts two places to the right (B->D, C->E).
$ws.Columns("B:C").Insert()

# Headers for the two freshly inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# No rating events have been recorded yet for these two new dates, so every
# data row just carries the same "UN" (unrated/no-change) placeholder that
# the rest of the sheet uses for an empty day.
$ws.Range("B2:C27").Value = "UN"

# Keep the same narrow, fixed width (~8 characters) used by the original
# date column across all three of the narrow data columns now in play.
$ws.Columns("C:E").ColumnWidth = 7.09
